$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 16, pushing the existing rows 16-28 down to 18-30.
$ws.Rows("16:17").Insert()

# New row 16: Argentina(o), 2021-10-15 (serial 44484)
$ws.Cells.Item(16,1).Value = 11
$ws.Cells.Item(16,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16,3).Value = "Bíobío"
$ws.Cells.Item(16,4).Value = 44484
$ws.Cells.Item(16,5).Value = 8
$ws.Cells.Item(16,6).Value = 100112013
$ws.Cells.Item(16,7).Value = "Alcachofa"
$ws.Cells.Item(16,8).Value = "Argentina(o)"
$ws.Cells.Item(16,9).Value = "Primera"
$ws.Cells.Item(16,10).Value = 220
$ws.Cells.Item(16,11).Value = 8000
$ws.Cells.Item(16,12).Value = 9000
$ws.Cells.Item(16,13).Value = 8455
$ws.Cells.Item(16,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(16,15).Value = "Provincia de Limarí"
$ws.Cells.Item(16,16).Value = 169
$ws.Cells.Item(16,17).Value = 50
$ws.Cells.Item(16,18).Value = "Hortaliza"

# New row 17: Española, 2021-10-15 (serial 44484)
$ws.Cells.Item(17,1).Value = 11
$ws.Cells.Item(17,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(17,3).Value = "Bíobío"
$ws.Cells.Item(17,4).Value = 44484
$ws.Cells.Item(17,5).Value = 8
$ws.Cells.Item(17,6).Value = 100112013
$ws.Cells.Item(17,7).Value = "Alcachofa"
$ws.Cells.Item(17,8).Value = "Española"
$ws.Cells.Item(17,9).Value = "Primera"
$ws.Cells.Item(17,10).Value = 220
$ws.Cells.Item(17,11).Value = 7500
$ws.Cells.Item(17,12).Value = 8000
$ws.Cells.Item(17,13).Value = 7727
$ws.Cells.Item(17,14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(17,15).Value = "Provincia de Limarí"
$ws.Cells.Item(17,16).Value = 258
$ws.Cells.Item(17,17).Value = 30
$ws.Cells.Item(17,18).Value = "Hortaliza"
